# Day-16/17, Completed BST sheet question: mark the remaining "Binary Search
# Trees" rows as Done, matching the already-completed rows above them.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Rows 214-218 & 220-222: mark "Done" AND pick up the green "Good"
# cell style (matching the already-completed rows just above, e.g. C202).
$doneAndStyledRows = @(214, 215, 216, 217, 218, 220, 221, 222)
foreach ($r in $doneAndStyledRows) {
    $ws.Range("C$r").Value = "Done"
}
$ws.Range("C202").Copy() | Out-Null
foreach ($r in $doneAndStyledRows) {
    $ws.Range("C$r").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}
$excel.CutCopyMode = $false

# Rows 223-233 & 235: just mark "Done", keep the existing (unstyled) format.
$doneRows = @(223, 224, 225, 226, 227, 228, 229, 230, 231, 232, 233, 235)
foreach ($r in $doneRows) {
    $ws.Range("C$r").Value = "Done"
}

# Scroll down to the newly completed section and select the next row to work on.
$excel.ActiveWindow.ScrollRow = 223
$ws.Range("B231").Select()
